$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''19.946.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -8.14%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''1.410.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -8.29%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = '''0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.04%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = '''  +0.01%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''274.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -5.23%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = '''0.3678'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -6.64%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.3124'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -2.41%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''39.36'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -9.48%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''1.037'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -3.78%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''0.06492'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -9.93%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''0.9997'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.05%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''5.458'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -5.41%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''17.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -4.46%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''6.184'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -6.86%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''1.409.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -8.25%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''0.00001015'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -7.48%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''0.05676'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -14.23%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''0.9997'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +0.29%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''70.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -16.34%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''5.599'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -9.03%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''14.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -5.24%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''11.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +2.13%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''2.281'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -3.67%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''19.946.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -8.17%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''2.259'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -5.65%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''135.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -10.54%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''16.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -8.36%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''1.565.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -8.46%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''109.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -7.22%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''4.089'
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = '''5.340'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -12.33%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''0.8207'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -15.47%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''0.07683'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -5.18%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''8.438'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.94%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = '''1.465'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -2.04%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''0.05831'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -2.43%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''4.871'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -6.51%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''1.001'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +0.06%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''0.02072'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -6.97%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''0.1910'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -6.87%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''10.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -7.64%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''1.089'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -8.00%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''0.5298'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -9.06%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''12.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -6.74%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''3.526'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -5.50%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''0.5158'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -7.81%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''113.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -2.58%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''1.771'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -6.47%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''1.040'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -10.61%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = '''1.002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.19%  '
$ws.Range("E51").Style = "Normal"
